$d = $word.ActiveDocument

# --- 1. Table indents (w:tblInd) -----------------------------------------
# Word's Rows.LeftIndent is expressed in points; the OOXML w:tblInd is in
# twips (dxa) = points * 20.
#   Table1 (Control de cambios):   40  dxa (2.0 pt)  -> -60 dxa (-3.0 pt)
#   Table2..Table5:                820 dxa (41.0 pt) -> 720 dxa (36.0 pt)
$d.Tables.Item(1).Rows.LeftIndent = -3
$d.Tables.Item(2).Rows.LeftIndent = 36
$d.Tables.Item(3).Rows.LeftIndent = 36
$d.Tables.Item(4).Rows.LeftIndent = 36
$d.Tables.Item(5).Rows.LeftIndent = 36

# --- 2. Alignment fixes inside the "Control de cambios" table (Table1) ---
# wdAlignParagraphCenter = 1
$t1 = $d.Tables.Item(1)
$t1.Cell(2, 2).Range.Paragraphs.Item(1).Alignment = 1   # 03/10/22
$t1.Cell(2, 3).Range.Paragraphs.Item(1).Alignment = 1   # Inicio
$t1.Cell(2, 4).Range.Paragraphs.Item(1).Alignment = 1   # Alberto Benítez
$t1.Cell(3, 2).Range.Paragraphs.Item(1).Alignment = 1   # 09/10/22
$t1.Cell(3, 3).Range.Paragraphs.Item(1).Alignment = 1   # Corrección
$t1.Cell(3, 4).Range.Paragraphs.Item(1).Alignment = 1   # Alberto Benítez

# --- 3. Text updates -------------------------------------------------------
$t3 = $d.Tables.Item(3)
$t3.Cell(2, 2).Range.Paragraphs.Item(1).Range.Text = "Plan de Dirección del Proyecto"
$t3.Cell(3, 2).Range.Paragraphs.Item(1).Range.Text = "Entrega final"

# --- 4. Remove the blank paragraph between "Productos relacionados
#        (sugerencias)" and the "Restricciones" heading -------------------
# NOTE: use $d.Content.Paragraphs (recomputed from the Content range) and
# not the cached $d.Paragraphs collection, which gets confused once a
# Tables member has been touched above.
$paras = $d.Content.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "Productos relacionados (sugerencias)`r") {
        $blank = $d.Content.Paragraphs.Item($i + 1)
        $blank.Range.Delete()
        break
    }
}
